$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("axes")

# Insert three new columns (D:F), shifting the existing "Title" column to G
$ws.Range("D1:F2").Insert(-4161)

# Fill in the new header/value columns
$ws.Range("D1").Value = "A_arrow"
$ws.Range("E1").Value = "B_arrow"
$ws.Range("F1").Value = "C_arrow"

$ws.Range("D2").Value = "Q (%)"
$ws.Range("E2").Value = "F (%)"
$ws.Range("F2").Value = "L (%)"

# Set the column widths for the new columns (Excel stores this as width 14.109375)
$ws.Range("D:F").ColumnWidth = 13.36

# Update selection and activate the axes sheet so it becomes the selected tab
$ws.Range("F4").Select()
$ws.Activate()
